$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues constant used for PasteSpecial below.
$xlPasteValues = -4163

function Set-TextCell {
    param($addr, $val)
    $escaped = $val -replace '"', '""'
    $range = $ws.Range($addr)
    # Quote the literal in a formula so Excel keeps it as text (the
    # sheet stores plain numeric-looking strings such as "1.00" or
    # "20.49" as text, not numbers), then convert the formula to a
    # static value via copy / paste-special so no formula or extra
    # cell-style is left behind.
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial($xlPasteValues) | Out-Null
}

$ws.Range('D2').Value = '76.341.66'
$ws.Range('E2').Value = '  +0.73%  '
$ws.Range('D3').Value = '3.036.10'
$ws.Range('E3').Value = '  +3.95%  '
$ws.Range('E4').Value = '  -0.05%  '
Set-TextCell 'D5' '198.10'
$ws.Range('E5').Value = '  -1.04%  '
Set-TextCell 'D6' '617.72'
$ws.Range('E6').Value = '  +3.46%  '
Set-TextCell 'D7' '0.999'
$ws.Range('E7').Value = '  +0.01%  '
Set-TextCell 'D8' '0.546'
$ws.Range('E8').Value = '  -1.06%  '
$ws.Range('E9').Value = '  +4.68%  '
$ws.Range('D10').Value = '3.032.90'
$ws.Range('E10').Value = '  +3.82%  '
Set-TextCell 'D11' '0.435'
$ws.Range('E11').Value = '  -1.58%  '
$ws.Range('E12').Value = '  -0.65%  '
Set-TextCell 'D13' '5.26'
$ws.Range('E13').Value = '  +6.79%  '
$ws.Range('D14').Value = '3.589.29'
$ws.Range('E14').Value = '  +3.76%  '
Set-TextCell 'D15' '28.76'
$ws.Range('E15').Value = '  +2.30%  '
$ws.Range('D16').Value = '76.306.66'
$ws.Range('E16').Value = '  +0.78%  '
$ws.Range('E17').Value = '  +1.98%  '
$ws.Range('D18').Value = '3.032.34'
$ws.Range('E18').Value = '  +3.79%  '
Set-TextCell 'D19' '13.46'
$ws.Range('E19').Value = '  +2.29%  '
Set-TextCell 'D20' '8.94'
$ws.Range('E20').Value = '  +2.47%  '
Set-TextCell 'D21' '379.40'
$ws.Range('E21').Value = '  +1.77%  '
$ws.Range('E22').Value = '  +2.06%  '
$ws.Range('E23').Value = '  +0.37%  '
$ws.Range('E24').Value = '  +3.31%  '
Set-TextCell 'D25' '72.74'
$ws.Range('E25').Value = '  +1.25%  '
Set-TextCell 'D26' '0.999'
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('E27').Value = '  +0.77%  '
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('E29').Value = '  +0.27%  '
Set-TextCell 'D30' '1.00'
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('E31').Value = '  +5.15%  '
$ws.Range('E32').Value = '  +0.79%  '
Set-TextCell 'D33' '489.79'
$ws.Range('E33').Value = '  -1.91%  '
Set-TextCell 'D34' '1.92'
$ws.Range('E34').Value = '  +4.15%  '
$ws.Range('E35').Value = '  -0.04%  '
Set-TextCell 'D36' '20.49'
$ws.Range('E36').Value = '  +1.52%  '
Set-TextCell 'D37' '161.85'
$ws.Range('E37').Value = '  -1.30%  '
$ws.Range('E39').Value = '  +4.65%  '
Set-TextCell 'D40' '0.381'
$ws.Range('E40').Value = '  +3.15%  '
Set-TextCell 'D41' '190.18'
$ws.Range('E41').Value = '  +6.79%  '
$ws.Range('E42').Value = '  -4.33%  '
$ws.Range('E43').Value = '  +0.01%  '
Set-TextCell 'D44' '0.787'
$ws.Range('E44').Value = '  +19.76%  '
Set-TextCell 'D45' '5.06'
$ws.Range('E45').Value = '  +2.14%  '
Set-TextCell 'D46' '41.95'
$ws.Range('E46').Value = '  +4.43%  '
$ws.Range('E47').Value = '  +4.98%  '
$ws.Range('E48').Value = '  -1.39%  '
$ws.Range('E49').Value = '  +3.84%  '
$ws.Range('E50').Value = '  +3.87%  '
$ws.Range('E51').Value = '  +1.02%  '

$excel.CutCopyMode = 0
